$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '65.936.85'
Set-TextValue $ws.Range("E2") '  -3.08%  '

Set-TextValue $ws.Range("D3") '3.486.16'
Set-TextValue $ws.Range("E3") '  +0.10%  '

Set-TextValue $ws.Range("E4") '  +0.02%  '

Set-TextValue $ws.Range("D5") '582.81'
Set-TextValue $ws.Range("E5") '  -1.91%  '

Set-TextValue $ws.Range("D6") '173.90'
Set-TextValue $ws.Range("E6") '  -4.74%  '

Set-TextValue $ws.Range("E7") '  +0.00%  '

Set-TextValue $ws.Range("D8") '0.599'
Set-TextValue $ws.Range("E8") '  -2.75%  '

Set-TextValue $ws.Range("D9") '3.485.34'
Set-TextValue $ws.Range("E9") '  +0.15%  '

Set-TextValue $ws.Range("E10") '  -7.10%  '

Set-TextValue $ws.Range("E11") '  -2.27%  '

Set-TextValue $ws.Range("D12") '0.411'
Set-TextValue $ws.Range("E12") '  -4.49%  '

Set-TextValue $ws.Range("D13") '4.086.42'
Set-TextValue $ws.Range("E13") '  +0.03%  '

Set-TextValue $ws.Range("E14") '  +0.17%  '

Set-TextValue $ws.Range("D15") '30.12'
Set-TextValue $ws.Range("E15") '  -6.69%  '

Set-TextValue $ws.Range("D16") '66.082.08'
Set-TextValue $ws.Range("E16") '  -2.73%  '

Set-TextValue $ws.Range("E17") '  -3.41%  '

Set-TextValue $ws.Range("D18") '3.479.88'
Set-TextValue $ws.Range("E18") '  -0.18%  '

Set-TextValue $ws.Range("D19") '5.95'
Set-TextValue $ws.Range("E19") '  -4.40%  '

Set-TextValue $ws.Range("D20") '13.93'
Set-TextValue $ws.Range("E20") '  -1.38%  '

Set-TextValue $ws.Range("D21") '366.27'
Set-TextValue $ws.Range("E21") '  -7.33%  '

Set-TextValue $ws.Range("D22") '7.77'
Set-TextValue $ws.Range("E22") '  -2.26%  '

Set-TextValue $ws.Range("D23") '72.76'
Set-TextValue $ws.Range("E23") '  +0.82%  '

Set-TextValue $ws.Range("E24") '  +0.15%  '

Set-TextValue $ws.Range("D25") '0.536'
Set-TextValue $ws.Range("E25") '  -0.88%  '

Set-TextValue $ws.Range("D26") '0.0000124'
Set-TextValue $ws.Range("E26") '  +1.54%  '

Set-TextValue $ws.Range("E27") '  -7.31%  '

Set-TextValue $ws.Range("E28") '  +0.76%  '

Set-TextValue $ws.Range("D29") '0.998'
Set-TextValue $ws.Range("E29") '  -0.28%  '

Set-TextValue $ws.Range("D30") '24.17'
Set-TextValue $ws.Range("E30") '  +2.36%  '

Set-TextValue $ws.Range("E31") '  -5.85%  '

Set-TextValue $ws.Range("E32") '  -3.83%  '

Set-TextValue $ws.Range("E33") '  +0.05%  '

Set-TextValue $ws.Range("D34") '7.15'
Set-TextValue $ws.Range("E34") '  -3.02%  '

Set-TextValue $ws.Range("E35") '  -8.29%  '

Set-TextValue $ws.Range("E36") '  -2.23%  '

Set-TextValue $ws.Range("D37") '160.12'
Set-TextValue $ws.Range("E37") '  -0.99%  '

Set-TextValue $ws.Range("D38") '29.34'
Set-TextValue $ws.Range("E38") '  +11.79%  '

Set-TextValue $ws.Range("D39") '0.889'
Set-TextValue $ws.Range("E39") '  -0.26%  '

Set-TextValue $ws.Range("D40") '2.821.49'
Set-TextValue $ws.Range("E40") '  +2.60%  '

Set-TextValue $ws.Range("E41") '  -5.83%  '

Set-TextValue $ws.Range("D42") '2.57'
Set-TextValue $ws.Range("E42") '  -9.09%  '

Set-TextValue $ws.Range("E43") '  -4.71%  '

Set-TextValue $ws.Range("D44") '6.41'
Set-TextValue $ws.Range("E44") '  -5.23%  '

Set-TextValue $ws.Range("E45") '  -4.97%  '

Set-TextValue $ws.Range("E46") '  -3.86%  '

Set-TextValue $ws.Range("E47") '  -8.11%  '

Set-TextValue $ws.Range("E48") '  -3.80%  '

Set-TextValue $ws.Range("D49") '307.03'
Set-TextValue $ws.Range("E49") '  -7.46%  '

Set-TextValue $ws.Range("E50") '  -3.35%  '

# Row 51: Stellar -> Cosmos (coin identity change)
Set-TextValue $ws.Range("B51") 'Cosmos'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D51") '6.22'
Set-TextValue $ws.Range("E51") '  -2.33%  '
